$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sprint 2")
$ws2 = $wb.Worksheets.Item("Sprint1")

# --- Update backlog items on "Sprint 2" (rows 24-26) to reflect the new sprint ---
$ws1.Range("A24").Value2 = "Home"
$ws1.Range("C24").Value2 = "Change confirmation windows to improve userability"

$ws1.Range("A25").Value2 = "Home"
$ws1.Range("C25").Value2 = "Implement Consistent Sidebar Navigation"

$ws1.Range("A26").Value2 = "Home"
$ws1.Range("C26").Value2 = "Implement employee summary across all pages"

# --- Remove the now-obsolete backlog rows (old rows 27-33) ---
$ws1.Rows("27:33").Delete()

# The H column total keeps referencing the old (pre-edit) range end, matching the source workbook
$ws1.Range("H27").Formula = "=SUM(H3:H24)"

# --- Update the burndown chart series to follow the Totals row to its new location ---
$co = $ws1.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,,'Sprint 2'!`$D`$27:`$H`$27,1)"

# --- Update selections to match the saved view state ---
$ws2.Activate()
$ws2.Range("M24").Select()
$ws1.Activate()
$ws1.Range("C27").Select()
